$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending/target cluster labels per data row (2..10), mirrored in cols A (sending) and D (target)
$clusters = @("ECs", "FAPs", "sCs")

# Full 3x3 grid of numeric metrics, keyed by [sending-cluster-index, target-cluster-index]
# Columns: E F G H I J K L M N O P Q R S T
$data = @{
  "0,0" = @(3,1,21.84976866666667,65.549306,0.05020018890879543,0.05020018890879543,2,0.6666666666666666,1.321445333333333,3.964336,0.01021782062667047,0.01021782062667047,28.87327483897955,259.859473550816,0.000512936525695044,0.0005129365256950441)
  "0,1" = @(3,1,21.84976866666667,65.549306,0.05020018890879543,0.05020018890879543,3,1,105.9632263333333,317.889679,0.819340166699254,0.8193401666992541,2315.271982556975,20837.44784301277,0.04113103114886649,0.04113103114886649)
  "0,2" = @(3,1,21.84976866666667,65.549306,0.05020018890879543,0.05020018890879543,3,1,22.04284166666666,66.128525,0.1704420126740755,0.1704420126740755,481.6309911726277,4334.67892055365,0.008556221234233894,0.008556221234233894)
  "1,0" = @(3,1,385.0524703333334,1155.157411,0.8846641374295412,0.8846641374295412,2,0.6666666666666666,1.321445333333333,3.964336,0.01021782062667047,0.01021782062667047,508.8257900104551,4579.432110094097,0.009039339471103204,0.009039339471103205)
  "1,1" = @(3,1,385.0524703333334,1155.157411,0.8846641374295412,0.8846641374295412,3,1,105.9632263333333,317.889679,0.819340166699254,0.8193401666992541,40801.40206414012,367212.6185772611,0.724840861834372,0.7248408618343721)
  "1,2" = @(3,1,385.0524703333334,1155.157411,0.8846641374295412,0.8846641374295412,3,1,22.04284166666666,66.128525,0.1704420126740755,0.1704420126740755,8487.650636916531,76388.85573224878,0.1507839361240659,0.1507839361240659)
  "2,0" = @(3,1,28.350479,85.05143699999999,0.06513567366166337,0.06513567366166337,2,0.6666666666666666,1.321445333333333,3.964336,0.01021782062667047,0.01021782062667047,37.46360817231466,337.172473550832,0.0006655446298722203,0.0006655446298722204)
  "2,1" = @(3,1,28.350479,85.05143699999999,0.06513567366166337,0.06513567366166337,3,1,105.9632263333333,317.889679,0.819340166699254,0.8193401666992541,3004.108222935413,27036.97400641872,0.05336827371601548,0.05336827371601548)
  "2,2" = @(3,1,28.350479,85.05143699999999,0.06513567366166337,0.06513567366166337,3,1,22.04284166666666,66.128525,0.1704420126740755,0.1704420126740755,624.9251197711582,5624.326077940425,0.01110185531577567,0.01110185531577567)
}

$numCols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$row = 2
for ($i = 0; $i -lt 3; $i++) {
  for ($j = 0; $j -lt 3; $j++) {
    $ws.Range("A$row").Value = $clusters[$i]
    $ws.Range("B$row").Value = "Fn1"
    $ws.Range("C$row").Value = "Sdc2"
    $ws.Range("D$row").Value = $clusters[$j]

    $vals = $data["$i,$j"]
    for ($k = 0; $k -lt $numCols.Length; $k++) {
      $ws.Range($numCols[$k] + "$row").Value = $vals[$k]
    }

    $row++
  }
}
